# Append a new data row (row 53) to each of the 4 worksheets, mirroring the
# structure/style of the existing last row (row 52), with the new record's
# values.

$wb = $excel.ActiveWorkbook

$dateVal = [double]"45839.49532407407"

$sheetsData = @(
    @{
        Index = 1
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x58"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 344
        I = 15
    },
    @{
        Index = 2
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x68"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 360
        I = 14
    },
    @{
        Index = 3
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x68"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 104
        I = 3
    },
    @{
        Index = 4
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x68"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 104
        I = 3
    }
)

foreach ($entry in $sheetsData) {
    $ws = $wb.Worksheets.Item($entry.Index)

    # Column A: date/time value, matching the number format used by row 52.
    $ws.Range("A53").Value = $dateVal
    $ws.Range("A53").NumberFormat = $ws.Range("A52").NumberFormat

    # Columns B-E: hex-string text fields.
    $ws.Range("B53").Value = $entry.B
    $ws.Range("C53").Value = $entry.C
    $ws.Range("D53").Value = $entry.D
    $ws.Range("E53").Value = $entry.E

    # Columns F-I: numeric fields.
    $ws.Range("F53").Value = $entry.F
    $ws.Range("G53").Value = $entry.G
    $ws.Range("H53").Value = $entry.H
    $ws.Range("I53").Value = $entry.I
}
